$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing row 2 with new test-candidate values ----
$ws.Range("A2").Value = "omLKn912"
$ws.Range("B2").Value = 2012454353
$ws.Range("C2").Value = "ngayxqc48"
$ws.Range("D2").Value = "Rf&N5n!8"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "UjDHiFTq"
$ws.Range("G2").Value = "ezdS"
$ws.Range("H2").Value = "Candidate"

# ---- Append new row 3 ----
$ws.Range("A3").Value = "sEdlG343"
$ws.Range("B3").Value = 2012454354
$ws.Range("C3").Value = "uymwjkw37"
$ws.Range("D3").Value = "hm65%E#P"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "xZVPOvsr"
$ws.Range("G3").Value = "WvJv"
$ws.Range("H3").Value = "Candidate"

# ---- Append new row 4 ----
$ws.Range("A4").Value = "iUkOr451"
$ws.Range("B4").Value = 2012454355
$ws.Range("C4").Value = "sghlmrc51"
$ws.Range("D4").Value = "xF&N3$2k"
$ws.Range("E4").Value = "MR"
$ws.Range("F4").Value = "ffHADbov"
$ws.Range("G4").Value = "fpzh"
$ws.Range("H4").Value = "Candidate"

# ---- Match the formatting used by the rest of the data rows (thin black
# borders, regular non-bold 11pt Calibri font) on the two newly added rows ----
$newRows = $ws.Range("A3:H4")
$newRows.Font.Name = "Calibri"
$newRows.Font.Size = 11
$newRows.Font.Bold = $false
$newRows.Font.Italic = $false
$newRows.Borders.LineStyle = 1
$newRows.Borders.Color = 0

# ---- Grow the visible selection to cover the whole (now larger) table ----
$ws.Range("A1:H4").Select() | Out-Null
